$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("D3").Value = "25 Oktober 2021"
$ws.Range("F3").Value = "DONE"
$ws.Range("E18").Select() | Out-Null
